$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) & Volume(1h) (E) figures to the latest scraped values.
# Rows 27/28 also swap which coin - EthereumClassic vs Stellar - occupies which rank.
# A leading apostrophe forces plain-decimal prices to stay text cells (matching the
# original inline-string formatting) instead of being auto-parsed as numbers.
$ws.Range("D2").Value = '27.930.36'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '1.637.02'
$ws.Range("E3").Value = '  -0.69%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'212.48"
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = "'23.31"
$ws.Range("E8").Value = '  -1.19%  '
$ws.Range("D9").Value = "'0.259"
$ws.Range("E9").Value = '  -2.40%  '
$ws.Range("E10").Value = '  -0.23%  '
$ws.Range("D11").Value = "'0.0881"
$ws.Range("E11").Value = '  +1.01%  '
$ws.Range("D12").Value = '1.869.22'
$ws.Range("E12").Value = '  -0.71%  '
$ws.Range("D13").Value = '1.637.28'
$ws.Range("E13").Value = '  -0.79%  '
$ws.Range("E14").Value = '  -0.30%  '
$ws.Range("D15").Value = "'0.568"
$ws.Range("E15").Value = '  +0.92%  '
$ws.Range("D16").Value = "'65.25"
$ws.Range("E16").Value = '  -0.73%  '
$ws.Range("D17").Value = '27.938.27'
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").Value = "'231.07"
$ws.Range("E18").Value = '  -0.42%  '
$ws.Range("E19").Value = '  -0.32%  '
$ws.Range("E20").Value = '  -1.94%  '
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").Value = "'10.41"
$ws.Range("E22").Value = '  -2.82%  '
$ws.Range("E23").Value = '  -0.44%  '
$ws.Range("E24").Value = '  -3.46%  '
$ws.Range("D25").Value = "'153.84"
$ws.Range("E25").Value = '  +1.19%  '
$ws.Range("E26").Value = '  +0.71%  '
$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").Value = "'0.111"
$ws.Range("E27").Value = '  -0.36%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = "'15.64"
$ws.Range("E28").Value = '  -0.53%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("E30").Value = '  -0.35%  '
$ws.Range("E31").Value = '  -0.74%  '
$ws.Range("E32").Value = '  +1.13%  '
$ws.Range("D33").Value = '1.408.95'
$ws.Range("E33").Value = '  -3.12%  '
$ws.Range("E34").Value = '  -1.46%  '
$ws.Range("E35").Value = '  +1.44%  '
$ws.Range("E36").Value = '  +1.52%  '
$ws.Range("D37").Value = "'0.970"
$ws.Range("E37").Value = '  +5.85%  '
$ws.Range("E38").Value = '  +0.40%  '
$ws.Range("D39").Value = "'0.562"
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("E40").Value = '  -1.66%  '
$ws.Range("E41").Value = '  +0.26%  '
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("D43").Value = "'67.12"
$ws.Range("E43").Value = '  -3.35%  '
$ws.Range("E44").Value = '  +2.56%  '
$ws.Range("E45").Value = '  +1.95%  '
$ws.Range("E46").Value = '  -1.82%  '
$ws.Range("D47").Value = '1.778.54'
$ws.Range("E47").Value = '  -0.74%  '
$ws.Range("D48").Value = "'87.97"
$ws.Range("E48").Value = '  -1.19%  '
$ws.Range("E49").Value = '  -0.27%  '
$ws.Range("E50").Value = '  -0.34%  '
$ws.Range("E51").Value = '  -1.63%  '
